$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update term labels in column D (rows 5-13 shift due to season
#     relabeling: seasonWinter removed, seasonSpring inserted) ---
$ws.Range("D5").Value = "seasonSpring"
$ws.Range("D6").Value = "seasonSummer"
$ws.Range("D7").Value = "seasonFall"
$ws.Range("D8").Value = "fish_basinWest:seasonSpring"
$ws.Range("D9").Value = "fish_basinNorth:seasonSpring"
$ws.Range("D10").Value = "fish_basinWest:seasonSummer"
$ws.Range("D11").Value = "fish_basinNorth:seasonSummer"
$ws.Range("D12").Value = "fish_basinWest:seasonFall"
$ws.Range("D13").Value = "fish_basinNorth:seasonFall"

# --- Update numeric results (estimate, std.error, statistic, p.value) ---

# Row 2 - (Intercept)
$ws.Range("E2").Value = 63.9490421663899
$ws.Range("F2").Value = 5.84634808873611
$ws.Range("G2").Value = 10.9382885171681
$ws.Range("H2").Value = 0.000000000000000000000000000756138161659184

# Row 3 - fish_basinWest
$ws.Range("E3").Value = -6.92482794038639
$ws.Range("F3").Value = 9.19983626300987
$ws.Range("G3").Value = -0.752712085564969
$ws.Range("H3").Value = 0.451622945903206

# Row 4 - fish_basinNorth
$ws.Range("E4").Value = 1.96775395662565
$ws.Range("F4").Value = 9.78924461576805
$ws.Range("G4").Value = 0.201011828170694
$ws.Range("H4").Value = 0.840689325293592

# Row 5 - seasonSpring
$ws.Range("E5").Value = 5.15522888053912
$ws.Range("F5").Value = 0.939241157249013
$ws.Range("G5").Value = 5.48871697194202
$ws.Range("H5").Value = 0.0000000404863738268761

# Row 6 - seasonSummer
$ws.Range("E6").Value = 18.9960979832281
$ws.Range("F6").Value = 0.984585974459425
$ws.Range("G6").Value = 19.2934883047239
$ws.Range("H6").Value = 0.0000000000000000000000000000000000000000000000000000000000000000000000000000000000609208018769435

# Row 7 - seasonFall
$ws.Range("E7").Value = 17.8745995915453
$ws.Range("F7").Value = 0.970139778635657
$ws.Range("G7").Value = 18.4247672193
$ws.Range("H7").Value = 0.000000000000000000000000000000000000000000000000000000000000000000000000000831475795147886

# Row 8 - fish_basinWest:seasonSpring
$ws.Range("E8").Value = 6.25981986905563
$ws.Range("F8").Value = 2.26548333187795
$ws.Range("G8").Value = 2.76312775334639
$ws.Range("H8").Value = 0.00572503536149729

# Row 9 - fish_basinNorth:seasonSpring
$ws.Range("E9").Value = 5.00009918075326
$ws.Range("F9").Value = 2.22682472613375
$ws.Range("G9").Value = 2.24539413545786
$ws.Range("H9").Value = 0.024742841569497

# Row 10 - fish_basinWest:seasonSummer
$ws.Range("E10").Value = 5.68592769756434
$ws.Range("F10").Value = 2.1782608863092
$ws.Range("G10").Value = 2.6103061085573
$ws.Range("H10").Value = 0.00904612360291998

# Row 11 - fish_basinNorth:seasonSummer
$ws.Range("E11").Value = 3.7172719328726
$ws.Range("F11").Value = 2.09387817856051
$ws.Range("G11").Value = 1.77530477700863
$ws.Range("H11").Value = 0.0758475770777054

# Row 12 - fish_basinWest:seasonFall
$ws.Range("E12").Value = 12.4692733429239
$ws.Range("F12").Value = 2.17351279653477
$ws.Range("G12").Value = 5.73692198306982
$ws.Range("H12").Value = 0.00000000964126405643401

# Row 13 - fish_basinNorth:seasonFall
$ws.Range("E13").Value = -6.70289405815238
$ws.Range("F13").Value = 2.10879038006366
$ws.Range("G13").Value = -3.17854923918519
$ws.Range("H13").Value = 0.00148014075630594

# Row 14 - sd__(Intercept) (estimate only; F/G/H remain empty)
$ws.Range("E14").Value = 15.2273995412183

# Row 15 - Residual (estimate only; F/G/H remain empty)
$ws.Range("E15").Value = 14.2527016998416
